$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics and p-values (column C = DM_Stat, column D = P_Value)
$ws.Range("C2").Value = -2.804576371848287
$ws.Range("D2").Value = 0.008269779080174233

$ws.Range("C3").Value = -1.265911391011916
$ws.Range("D3").Value = 0.2141504887976426

$ws.Range("C4").Value = -1.062767487446342
$ws.Range("D4").Value = 0.2953779265888301

$ws.Range("C5").Value = -0.1008890171654989
$ws.Range("D5").Value = 0.9202311652240653

$ws.Range("C6").Value = 1.604527576742994
$ws.Range("D6").Value = 0.1178478837544774

$ws.Range("C7").Value = 1.597934885378422
$ws.Range("D7").Value = 0.119310952010006

$ws.Range("C8").Value = 2.910323198002666
$ws.Range("D8").Value = 0.006327709065648612

$ws.Range("C9").Value = 0.2323940787152471
$ws.Range("D9").Value = 0.8176251744173162

$ws.Range("C10").Value = 0.8388666675749418
$ws.Range("D10").Value = 0.407403519460789

$ws.Range("C11").Value = 0.6966178590831357
$ws.Range("D11").Value = 0.4907775309287745
